$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 685.3782603333333
$ws.Range("H2").Value = 2056.134781
$ws.Range("I2").Value = 0.7131285654702259
$ws.Range("J2").Value = 0.7131285654702259
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 17.16653
$ws.Range("N2").Value = 51.49959
$ws.Range("O2").Value = 0.0560345397128279
$ws.Range("P2").Value = 0.0560345397128279
$ws.Range("Q2").Value = 11765.56646735997
$ws.Range("R2").Value = 105890.0982062398
$ws.Range("S2").Value = 0.03995983092219336
$ws.Range("T2").Value = 0.03995983092219336
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 685.3782603333333
$ws.Range("H3").Value = 2056.134781
$ws.Range("I3").Value = 0.7131285654702259
$ws.Range("J3").Value = 0.7131285654702259
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 256.4443053333333
$ws.Range("N3").Value = 769.332916
$ws.Range("O3").Value = 0.8370788162388805
$ws.Range("P3").Value = 0.8370788162388805
$ws.Range("Q3").Value = 175761.3518617501
$ws.Range("R3").Value = 1581852.166755751
$ws.Range("S3").Value = 0.5969448154099476
$ws.Range("T3").Value = 0.5969448154099476
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 685.3782603333333
$ws.Range("H4").Value = 2056.134781
$ws.Range("I4").Value = 0.7131285654702259
$ws.Range("J4").Value = 0.7131285654702259
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 32.74538866666666
$ws.Range("N4").Value = 98.236166
$ws.Range("O4").Value = 0.1068866440482915
$ws.Range("P4").Value = 0.1068866440482915
$ws.Range("Q4").Value = 22442.97751829885
$ws.Range("R4").Value = 201986.7976646896
$ws.Range("S4").Value = 0.07622391913808479
$ws.Range("T4").Value = 0.07622391913808479
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 213.8079683333333
$ws.Range("H5").Value = 641.423905
$ws.Range("I5").Value = 0.2224648468854243
$ws.Range("J5").Value = 0.2224648468854243
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 17.16653
$ws.Range("N5").Value = 51.49959
$ws.Range("O5").Value = 0.0560345397128279
$ws.Range("P5").Value = 0.0560345397128279
$ws.Range("Q5").Value = 3670.340902633216
$ws.Range("R5").Value = 33033.06812369895
$ws.Range("S5").Value = 0.01246571529750949
$ws.Range("T5").Value = 0.01246571529750949
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 213.8079683333333
$ws.Range("H6").Value = 641.423905
$ws.Range("I6").Value = 0.2224648468854243
$ws.Range("J6").Value = 0.2224648468854243
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 256.4443053333333
$ws.Range("N6").Value = 769.332916
$ws.Range("O6").Value = 0.8370788162388805
$ws.Range("P6").Value = 0.8370788162388805
$ws.Range("Q6").Value = 54829.83591397299
$ws.Range("R6").Value = 493468.5232257569
$ws.Range("S6").Value = 0.1862206106856148
$ws.Range("T6").Value = 0.1862206106856148
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 213.8079683333333
$ws.Range("H7").Value = 641.423905
$ws.Range("I7").Value = 0.2224648468854243
$ws.Range("J7").Value = 0.2224648468854243
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 32.74538866666666
$ws.Range("N7").Value = 98.236166
$ws.Range("O7").Value = 0.1068866440482915
$ws.Range("P7").Value = 0.1068866440482915
$ws.Range("Q7").Value = 7001.225023105358
$ws.Range("R7").Value = 63011.02520794823
$ws.Range("S7").Value = 0.02377852090230002
$ws.Range("T7").Value = 0.02377852090230002
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 61.90030400000001
$ws.Range("H8").Value = 185.700912
$ws.Range("I8").Value = 0.06440658764434989
$ws.Range("J8").Value = 0.06440658764434989
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 17.16653
$ws.Range("N8").Value = 51.49959
$ws.Range("O8").Value = 0.0560345397128279
$ws.Range("P8").Value = 0.0560345397128279
$ws.Range("Q8").Value = 1062.61342562512
$ws.Range("R8").Value = 9563.520830626081
$ws.Range("S8").Value = 0.003608993493125055
$ws.Range("T8").Value = 0.003608993493125055
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 61.90030400000001
$ws.Range("H9").Value = 185.700912
$ws.Range("I9").Value = 0.06440658764434989
$ws.Range("J9").Value = 0.06440658764434989
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 256.4443053333333
$ws.Range("N9").Value = 769.332916
$ws.Range("O9").Value = 0.8370788162388805
$ws.Range("P9").Value = 0.8370788162388805
$ws.Range("Q9").Value = 15873.98045920216
$ws.Range("R9").Value = 142865.8241328194
$ws.Range("S9").Value = 0.05391339014331811
$ws.Range("T9").Value = 0.05391339014331811
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 61.90030400000001
$ws.Range("H10").Value = 185.700912
$ws.Range("I10").Value = 0.06440658764434989
$ws.Range("J10").Value = 0.06440658764434989
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 32.74538866666666
$ws.Range("N10").Value = 98.236166
$ws.Range("O10").Value = 0.1068866440482915
$ws.Range("P10").Value = 0.1068866440482915
$ws.Range("Q10").Value = 2026.949513064821
$ws.Range("R10").Value = 18242.54561758339
$ws.Range("S10").Value = 0.006884204007906718
$ws.Range("T10").Value = 0.006884204007906718
